$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13: Inscritos, Pagos, Inscrições homologadas +1
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 5
$ws.Range("H13").Value = 5

# Row 14: Pagos, Inscrições homologadas -1
$ws.Range("F14").Value = 16
$ws.Range("H14").Value = 19

# Row 15: Inscritos, Pagos, Inscrições homologadas +1
$ws.Range("E15").Value = 93
$ws.Range("F15").Value = 45
$ws.Range("H15").Value = 56
